$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Piña" dataset. It belongs
# chronologically before the existing row 384, so insert a new row there
# and push the rest of the table (old rows 384-396) down by one, then
# populate the newly inserted row with its own data.
$ws.Rows.Item(384).Insert()

$ws.Range("A384").Value = 5
$ws.Range("B384").Value = "Macroferia Regional de Talca"
$ws.Range("C384").Value = "Maule"
$ws.Range("D384").Value = 45075
$ws.Range("E384").Value = 7
$ws.Range("F384").Value = "Fruta"
$ws.Range("G384").Value = 100108
$ws.Range("H384").Value = "Tropicales y subtropicales"
$ws.Range("I384").Value = 100108005
$ws.Range("J384").Value = "Piña"
$ws.Range("K384").Value = "Caramelo"
$ws.Range("L384").Value = "Segunda"
$ws.Range("M384").Value = 200
$ws.Range("N384").Value = 13000
$ws.Range("O384").Value = 13000
$ws.Range("P384").Value = 13000
$ws.Range("Q384").Value = "`$/caja 14 unidades"
$ws.Range("R384").Value = "Ecuador"
$ws.Range("S384").Value = 929
$ws.Range("T384").Value = 14
